$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the existing header cell (G1) onto the new header cell (H1)
# so the new "Save" header reuses the same bold/border/centered style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for the new "Save" column, rows 2 through 13
$values = @(1, 0, 0, 0, 0, 1, 1, 1, 1, 0, 1, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
